$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4860750734806061
$ws.Range("B1").Value = 1.29093337059021
$ws.Range("C1").Value = 3.742969512939453
$ws.Range("D1").Value = 3.20445442199707
$ws.Range("E1").Value = 0.8148506879806519
